$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.197.39"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").Value = "3.524.43"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'608.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").Value = "'148.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.13%  "

$ws.Range("D7").Value = "3.524.09"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").Value = "'7.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.69%  "

$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "'0.0000218"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").Value = "4.121.12"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").Value = "'31.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("D16").Value = "3.520.63"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "67.149.99"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "'10.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.31%  "

$ws.Range("D20").Value = "'6.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "'15.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("D22").Value = "'439.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.57%  "

$ws.Range("E23").Value = "  -2.78%  "

$ws.Range("D24").Value = "'79.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").Value = "3.666.51"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  -3.88%  "

$ws.Range("D28").Value = "'9.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.21%  "

$ws.Range("D29").Value = "'8.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.79%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -3.19%  "

$ws.Range("E32").Value = "  -3.02%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'25.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").Value = "3.516.59"
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("E36").Value = "  -2.82%  "

$ws.Range("D37").Value = "'5.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.41%  "

$ws.Range("D38").Value = "'8.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.83%  "

$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").Value = "'0.0896"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D42").Value = "'172.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("D43").Value = "'5.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("E44").Value = "  -9.42%  "

$ws.Range("E45").Value = "  +0.50%  "

$ws.Range("E46").Value = "  -0.86%  "

$ws.Range("D47").Value = "'28.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.82%  "

$ws.Range("D48").Value = "'1.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.86%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.53%  "

$ws.Range("D51").Value = "'0.994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.38%  "
